# Novos filtros para gestao de leitos.
# Adds new "evaldo" grants (tb_crtr_intnc, tb_dieta, tb_const) and a full
# "GRANT SELECT on integracao.tb_const" block for every existing user on
# sheet "grants por usuario".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")
$ws.Activate()

# --- New table grants for "evaldo" (rows 81-83) -----------------------
# Column A text is written first in this particular order so the
# workbook's shared-strings table ends up with the same index order
# as the target file (crtr_intnc, dieta, administrativo, const, evaldo).
$ws.Range("A81").Value = "GRANT SELECT on integracao.tb_crtr_intnc TO "
$ws.Range("A82").Value = "GRANT SELECT on integracao.tb_dieta TO "
$ws.Range("B109").Value = "administrativo"
$ws.Range("A83").Value = "GRANT SELECT on integracao.tb_const TO "
$ws.Range("B81").Value = "evaldo"
$ws.Range("B82").Value = "evaldo"
$ws.Range("B83").Value = "evaldo"

$ws.Range("D81").Formula = '=A81&" "&B81&" "&C81'
$ws.Range("D82").Formula = '=A82&" "&B82&" "&C82'
$ws.Range("D83").Formula = '=A83&" "&B83&" "&C83'

# --- Remove the now-unused blank filler rows 84 and 85 -----------------
# Clearing every cell (incl. the leftover "=A&B&C" formula) leaves the
# rows completely empty, so they disappear from the saved worksheet.
$ws.Range("A84:D85").ClearContents()

# --- Row 86 stays as the existing blank filler row (A:C empty, D keeps
#     its "=A86&" "&B86&" "&C86" formula evaluating to "  ") -----------

# --- "GRANT SELECT on integracao.tb_const TO <user> ;" for every user --
$users = @(
    "aoliveira",
    "dalves",
    "emenezes",
    "gcassia",
    "lmaria",
    "mrezende",
    "lvieira",
    "tsilva",
    "vrodrigues",
    "vlucia",
    "vsilva",
    "woliveira",
    "wquetz",
    "ftesta",
    "simone",
    "grazielle",
    "dayane",
    "ronan",
    "clovismelo",
    "mariabethania",
    "fernandazeferino",
    "camila",
    "administrativo",
    "tivilaverde",
    "lamorim",
    "mvilela",
    "fcampos",
    "bcorrea",
    "mmattos",
    "greis",
    "ldelgado",
    "deliza",
    "aalbino",
    "ralmeida",
    "bsouza",
    "tnovaes",
    "mliberato",
    "dchinelato",
    "amonteiro",
    "soliveira"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = 87 + $i
    $ws.Range("A$row").Value = "GRANT SELECT on integracao.tb_const TO "
    $ws.Range("B$row").Value = $users[$i]
    $ws.Range("C$row").Value = ";"
    $ws.Range("D$row").Formula = '=A' + $row + '&" "&B' + $row + '&" "&C' + $row
}

# --- View state: scrolled so row 10 is at the top, C27 selected --------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
